$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value = 8300  # H32 (was 8666.666999999999)
$ws.Cells.Item(32, 9).Value = 7750.5  # I32 (was 8001)
$ws.Cells.Item(32, 10).Value = 8666.333000000001  # J32 (was 8999.5)
$ws.Cells.Item(32, 11).Value = 7750.5  # K32 (was 8001)
$ws.Cells.Item(32, 12).Value = 8666.333000000001  # L32 (was 8999.5)
$ws.Cells.Item(32, 13).Value = -7424.5  # M32 (was -7675)
$ws.Cells.Item(32, 14).Value = -9318.333000000001  # N32 (was -9651.5)
$ws.Cells.Item(38, 8).Value = 652.73334  # H38 (was 618.1875)
$ws.Cells.Item(38, 9).Value = 243.75  # I38 (was 232.6923)
$ws.Cells.Item(38, 10).Value = 2288.6667  # J38 (was 2288.6667)
$ws.Cells.Item(38, 11).Value = 731.25  # K38 (was 698.0769)
$ws.Cells.Item(38, 12).Value = 6866.000100000001  # L38 (was 6866.000100000001)
$ws.Cells.Item(38, 13).Value = -359.25  # M38 (was -326.0769)
$ws.Cells.Item(38, 14).Value = -7610.000100000001  # N38 (was -7610.000100000001)
$ws.Cells.Item(48, 8).Value = 2998.5  # H48 (was 3000)
$ws.Cells.Item(48, 9).Value = 3000  # I48 (was 3000)
$ws.Cells.Item(48, 10).Value = 2997  # J48 (was 0)
$ws.Cells.Item(48, 11).Value = 9000  # K48 (was 9000)
$ws.Cells.Item(48, 12).Value = 8991  # L48 (was 0)
$ws.Cells.Item(48, 13).Value = -8708  # M48 (was -8708)
$ws.Cells.Item(48, 14).Value = -9575  # N48 (was None)
$ws.Cells.Item(51, 8).Value = 22566.6  # H51 (was 21468.625)
$ws.Cells.Item(51, 9).Value = 6700  # I51 (was 6214)
$ws.Cells.Item(51, 10).Value = 30499.9  # J51 (was 33333.332)
$ws.Cells.Item(51, 11).Value = 6700  # K51 (was 6214)
$ws.Cells.Item(51, 12).Value = 30499.9  # L51 (was 33333.332)
$ws.Cells.Item(51, 13).Value = -6216  # M51 (was -5730)
$ws.Cells.Item(51, 14).Value = -31467.9  # N51 (was -34301.332)
$ws.Cells.Item(56, 8).Value = 2998.5  # H56 (was 3000)
$ws.Cells.Item(56, 9).Value = 3000  # I56 (was 3000)
$ws.Cells.Item(56, 10).Value = 2997  # J56 (was 0)
$ws.Cells.Item(56, 11).Value = 9000  # K56 (was 9000)
$ws.Cells.Item(56, 12).Value = 8991  # L56 (was 0)
$ws.Cells.Item(56, 13).Value = -8466  # M56 (was -8466)
$ws.Cells.Item(56, 14).Value = -10059  # N56 (was None)
$ws.Cells.Item(86, 8).Value = 2847.7144  # H86 (was 2776.2144)
$ws.Cells.Item(86, 9).Value = 2116.3333  # I86 (was 1871)
$ws.Cells.Item(86, 10).Value = 3396.25  # J86 (was 3681.4285)
$ws.Cells.Item(86, 11).Value = 2116.3333  # K86 (was 1871)
$ws.Cells.Item(86, 12).Value = 3396.25  # L86 (was 3681.4285)
$ws.Cells.Item(86, 13).Value = -993.3332999999998  # M86 (was -748)
$ws.Cells.Item(86, 14).Value = -5642.25  # N86 (was -5927.4285)
$ws.Cells.Item(89, 8).Value = 2847.7144  # H89 (was 2776.2144)
$ws.Cells.Item(89, 9).Value = 2116.3333  # I89 (was 1871)
$ws.Cells.Item(89, 10).Value = 3396.25  # J89 (was 3681.4285)
$ws.Cells.Item(89, 11).Value = 10581.6665  # K89 (was 9355)
$ws.Cells.Item(89, 12).Value = 16981.25  # L89 (was 18407.1425)
$ws.Cells.Item(89, 13).Value = -4965.666499999999  # M89 (was -3739)
$ws.Cells.Item(89, 14).Value = -28213.25  # N89 (was -29639.1425)
$ws.Cells.Item(100, 8).Value = 8692.308000000001  # H100 (was 7100)
$ws.Cells.Item(100, 9).Value = 2748.5  # I100 (was 2748.5)
$ws.Cells.Item(100, 10).Value = 9371.6  # J100 (was 10001)
$ws.Cells.Item(100, 11).Value = 2748.5  # K100 (was 2748.5)
$ws.Cells.Item(100, 12).Value = 9371.6  # L100 (was 10001)
$ws.Cells.Item(100, 13).Value = -2207.5  # M100 (was -2207.5)
$ws.Cells.Item(100, 14).Value = -10453.6  # N100 (was -11083)
$ws.Cells.Item(132, 8).Value = 1347  # H132 (was 1407.0834)
$ws.Cells.Item(132, 9).Value = 1321  # I132 (was 1387.9)
$ws.Cells.Item(132, 10).Value = 1503  # J132 (was 1503)
$ws.Cells.Item(132, 11).Value = 3963  # K132 (was 4163.700000000001)
$ws.Cells.Item(132, 12).Value = 4509  # L132 (was 4509)
$ws.Cells.Item(132, 13).Value = -1433  # M132 (was -1633.700000000001)
$ws.Cells.Item(132, 14).Value = -9569  # N132 (was -9569)
$ws.Cells.Item(137, 8).Value = 4169069  # H137 (was 4239706)
$ws.Cells.Item(137, 9).Value = 5104031.5  # I137 (was 5321162.5)
$ws.Cells.Item(137, 10).Value = 4236.4546  # J137 (was 4000.0833)
$ws.Cells.Item(137, 11).Value = 15312094.5  # K137 (was 15963487.5)
$ws.Cells.Item(137, 12).Value = 12709.3638  # L137 (was 12000.2499)
$ws.Cells.Item(137, 13).Value = -15309544.5  # M137 (was -15960937.5)
$ws.Cells.Item(137, 14).Value = -17809.3638  # N137 (was -17100.2499)
$ws.Cells.Item(138, 8).Value = 1841.262  # H138 (was 2014.1555)
$ws.Cells.Item(138, 9).Value = 1388.625  # I138 (was 1493.9062)
$ws.Cells.Item(138, 10).Value = 3289.7  # J138 (was 3294.7693)
$ws.Cells.Item(138, 11).Value = 4165.875  # K138 (was 4481.7186)
$ws.Cells.Item(138, 12).Value = 9869.099999999999  # L138 (was 9884.3079)
$ws.Cells.Item(138, 13).Value = 974.125  # M138 (was 658.2813999999998)
$ws.Cells.Item(138, 14).Value = -20149.1  # N138 (was -20164.3079)

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 7078555  # H2 (was 6572958.5)
$ws.Cells.Item(2, 9).Value = 7078555  # I2 (was 6572958.5)
$ws.Cells.Item(2, 10).Value = 0  # J2 (was 0)
$ws.Cells.Item(2, 11).Value = 7078555  # K2 (was 6572958.5)
$ws.Cells.Item(2, 12).Value = 0  # L2 (was 0)
$ws.Cells.Item(2, 13).Value = -7078442  # M2 (was -6572845.5)
$ws.Cells.Item(5, 8).Value = 269  # H5 (was 258.35715)
$ws.Cells.Item(5, 9).Value = 279.7  # I5 (was 281.7)
$ws.Cells.Item(5, 10).Value = 233.33333  # J5 (was 200)
$ws.Cells.Item(5, 11).Value = 279.7  # K5 (was 281.7)
$ws.Cells.Item(5, 12).Value = 233.33333  # L5 (was 200)
$ws.Cells.Item(5, 13).Value = -167.7  # M5 (was -169.7)
$ws.Cells.Item(5, 14).Value = -457.33333  # N5 (was -424)
$ws.Cells.Item(45, 8).Value = 24127.117  # H45 (was 23880.5)
$ws.Cells.Item(45, 9).Value = 41041.777  # I45 (was 37066.4)
$ws.Cells.Item(45, 10).Value = 5098.125  # J45 (was 7398.125)
$ws.Cells.Item(45, 11).Value = 41041.777  # K45 (was 37066.4)
$ws.Cells.Item(45, 12).Value = 5098.125  # L45 (was 7398.125)
$ws.Cells.Item(45, 13).Value = -40664.777  # M45 (was -36689.4)
$ws.Cells.Item(45, 14).Value = -5852.125  # N45 (was -8152.125)
$ws.Cells.Item(116, 8).Value = 7078555  # H116 (was 6572958.5)
$ws.Cells.Item(116, 9).Value = 7078555  # I116 (was 6572958.5)
$ws.Cells.Item(116, 10).Value = 0  # J116 (was 0)
$ws.Cells.Item(116, 11).Value = 7078555  # K116 (was 6572958.5)
$ws.Cells.Item(116, 12).Value = 0  # L116 (was 0)
$ws.Cells.Item(116, 13).Value = -7076261  # M116 (was -6570664.5)

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 7078555  # H3 (was 6572958.5)
$ws.Cells.Item(3, 9).Value = 7078555  # I3 (was 6572958.5)
$ws.Cells.Item(3, 10).Value = 0  # J3 (was 0)
$ws.Cells.Item(3, 11).Value = 7078555  # K3 (was 6572958.5)
$ws.Cells.Item(3, 12).Value = 0  # L3 (was 0)
$ws.Cells.Item(3, 13).Value = -7078441  # M3 (was -6572844.5)
$ws.Cells.Item(4, 8).Value = 269  # H4 (was 258.35715)
$ws.Cells.Item(4, 9).Value = 279.7  # I4 (was 281.7)
$ws.Cells.Item(4, 10).Value = 233.33333  # J4 (was 200)
$ws.Cells.Item(4, 11).Value = 279.7  # K4 (was 281.7)
$ws.Cells.Item(4, 12).Value = 233.33333  # L4 (was 200)
$ws.Cells.Item(4, 13).Value = -164.7  # M4 (was -166.7)
$ws.Cells.Item(4, 14).Value = -463.33333  # N4 (was -430)

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 404.4524  # H7 (was 399.69766)
$ws.Cells.Item(7, 9).Value = 433.4516  # I7 (was 433.4516)
$ws.Cells.Item(7, 10).Value = 322.72726  # J7 (was 312.5)
$ws.Cells.Item(7, 11).Value = 433.4516  # K7 (was 433.4516)
$ws.Cells.Item(7, 12).Value = 322.72726  # L7 (was 312.5)
$ws.Cells.Item(7, 13).Value = -320.4516  # M7 (was -320.4516)
$ws.Cells.Item(7, 14).Value = -548.72726  # N7 (was -538.5)
$ws.Cells.Item(31, 8).Value = 22730498  # H31 (was 25003648)
$ws.Cells.Item(31, 9).Value = 34485184  # I31 (was 41669370)
$ws.Cells.Item(31, 10).Value = 4771.8  # J31 (was 5067.3125)
$ws.Cells.Item(31, 11).Value = 34485184  # K31 (was 41669370)
$ws.Cells.Item(31, 12).Value = 4771.8  # L31 (was 5067.3125)
$ws.Cells.Item(31, 13).Value = -34484889  # M31 (was -41669075)
$ws.Cells.Item(31, 14).Value = -5361.8  # N31 (was -5657.3125)
$ws.Cells.Item(34, 8).Value = 22730498  # H34 (was 25003648)
$ws.Cells.Item(34, 9).Value = 34485184  # I34 (was 41669370)
$ws.Cells.Item(34, 10).Value = 4771.8  # J34 (was 5067.3125)
$ws.Cells.Item(34, 11).Value = 34485184  # K34 (was 41669370)
$ws.Cells.Item(34, 12).Value = 4771.8  # L34 (was 5067.3125)
$ws.Cells.Item(34, 13).Value = -34484982  # M34 (was -41669168)
$ws.Cells.Item(34, 14).Value = -5175.8  # N34 (was -5471.3125)
$ws.Cells.Item(122, 8).Value = 18443  # H122 (was 24485.125)
$ws.Cells.Item(122, 9).Value = 20874.777  # I122 (was 30146.834)
$ws.Cells.Item(122, 10).Value = 7500  # J122 (was 7500)
$ws.Cells.Item(122, 11).Value = 62624.33099999999  # K122 (was 90440.50199999999)
$ws.Cells.Item(122, 12).Value = 22500  # L122 (was 22500)
$ws.Cells.Item(122, 13).Value = -60174.33099999999  # M122 (was -87990.50199999999)
$ws.Cells.Item(122, 14).Value = -27400  # N122 (was -27400)
$ws.Cells.Item(134, 8).Value = 3319.2126  # H134 (was 3386.0652)
$ws.Cells.Item(134, 9).Value = 2721.2163  # I134 (was 2790.0278)
$ws.Cells.Item(134, 10).Value = 5531.8  # J134 (was 5531.8)
$ws.Cells.Item(134, 11).Value = 8163.6489  # K134 (was 8370.0834)
$ws.Cells.Item(134, 12).Value = 16595.4  # L134 (was 16595.4)
$ws.Cells.Item(134, 13).Value = -5628.6489  # M134 (was -5835.0834)
$ws.Cells.Item(134, 14).Value = -21665.4  # N134 (was -21665.4)

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 8).Value = 768.4286  # H68 (was 767.8570999999999)
$ws.Cells.Item(68, 9).Value = 1036  # I68 (was 1036)
$ws.Cells.Item(68, 10).Value = 661.4  # J68 (was 660.6)
$ws.Cells.Item(68, 11).Value = 3108  # K68 (was 3108)
$ws.Cells.Item(68, 12).Value = 1984.2  # L68 (was 1981.8)
$ws.Cells.Item(68, 13).Value = -2297  # M68 (was -2297)
$ws.Cells.Item(68, 14).Value = -3606.2  # N68 (was -3603.8)
$ws.Cells.Item(71, 8).Value = 768.4286  # H71 (was 767.8570999999999)
$ws.Cells.Item(71, 9).Value = 1036  # I71 (was 1036)
$ws.Cells.Item(71, 10).Value = 661.4  # J71 (was 660.6)
$ws.Cells.Item(71, 11).Value = 9324  # K71 (was 9324)
$ws.Cells.Item(71, 12).Value = 5952.599999999999  # L71 (was 5945.400000000001)
$ws.Cells.Item(71, 13).Value = -5268  # M71 (was -5268)
$ws.Cells.Item(71, 14).Value = -14064.6  # N71 (was -14057.4)
$ws.Cells.Item(80, 8).Value = 3212.75  # H80 (was 3374.625)
$ws.Cells.Item(80, 9).Value = 3000  # I80 (was 2500)
$ws.Cells.Item(80, 10).Value = 3243.1428  # J80 (was 3666.1667)
$ws.Cells.Item(80, 11).Value = 9000  # K80 (was 7500)
$ws.Cells.Item(80, 12).Value = 9729.428400000001  # L80 (was 10998.5001)
$ws.Cells.Item(80, 13).Value = -8064  # M80 (was -6564)
$ws.Cells.Item(80, 14).Value = -11601.4284  # N80 (was -12870.5001)
$ws.Cells.Item(83, 8).Value = 3212.75  # H83 (was 3374.625)
$ws.Cells.Item(83, 9).Value = 3000  # I83 (was 2500)
$ws.Cells.Item(83, 10).Value = 3243.1428  # J83 (was 3666.1667)
$ws.Cells.Item(83, 11).Value = 27000  # K83 (was 22500)
$ws.Cells.Item(83, 12).Value = 29188.2852  # L83 (was 32995.5003)
$ws.Cells.Item(83, 13).Value = -22320  # M83 (was -17820)
$ws.Cells.Item(83, 14).Value = -38548.2852  # N83 (was -42355.5003)
$ws.Cells.Item(113, 8).Value = 1675.3928  # H113 (was 1718.963)
$ws.Cells.Item(113, 9).Value = 1992.2858  # I113 (was 1818.125)
$ws.Cells.Item(113, 10).Value = 1569.762  # J113 (was 1677.2106)
$ws.Cells.Item(113, 11).Value = 5976.857400000001  # K113 (was 5454.375)
$ws.Cells.Item(113, 12).Value = 4709.286  # L113 (was 5031.6318)
$ws.Cells.Item(113, 13).Value = -3806.857400000001  # M113 (was -3284.375)
$ws.Cells.Item(113, 14).Value = -9049.286  # N113 (was -9371.631799999999)
$ws.Cells.Item(137, 8).Value = 3328.8  # H137 (was 3368.9)
$ws.Cells.Item(137, 9).Value = 2449.5  # I137 (was 1900)
$ws.Cells.Item(137, 10).Value = 3548.625  # J137 (was 3736.125)
$ws.Cells.Item(137, 11).Value = 7348.5  # K137 (was 5700)
$ws.Cells.Item(137, 12).Value = 10645.875  # L137 (was 11208.375)
$ws.Cells.Item(137, 13).Value = -2248.5  # M137 (was -600)
$ws.Cells.Item(137, 14).Value = -20845.875  # N137 (was -21408.375)
$ws.Cells.Item(140, 8).Value = 1233.0769  # H140 (was 1147.1428)
$ws.Cells.Item(140, 9).Value = 1004.2857  # I140 (was 1006)
$ws.Cells.Item(140, 10).Value = 1500  # J140 (was 1500)
$ws.Cells.Item(140, 11).Value = 3012.8571  # K140 (was 3018)
$ws.Cells.Item(140, 12).Value = 4500  # L140 (was 4500)
$ws.Cells.Item(140, 13).Value = 2167.1429  # M140 (was 2162)
$ws.Cells.Item(140, 14).Value = -14860  # N140 (was -14860)
$ws.Cells.Item(141, 8).Value = 5681.875  # H141 (was 6113.375)
$ws.Cells.Item(141, 9).Value = 5707.857  # I141 (was 6201)
$ws.Cells.Item(141, 10).Value = 5500  # J141 (was 5500)
$ws.Cells.Item(141, 11).Value = 17123.571  # K141 (was 18603)
$ws.Cells.Item(141, 12).Value = 16500  # L141 (was 16500)
$ws.Cells.Item(141, 13).Value = -11943.571  # M141 (was -13423)
$ws.Cells.Item(141, 14).Value = -26860  # N141 (was -26860)

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(17, 8).Value = 7894.143  # H17 (was 6860.5835)
$ws.Cells.Item(17, 9).Value = 0  # I17 (was 599)
$ws.Cells.Item(17, 10).Value = 7894.143  # J17 (was 7429.8184)
$ws.Cells.Item(17, 11).Value = 0  # K17 (was 599)
$ws.Cells.Item(17, 12).ClearContents()  # L17 (was 7429.8184)
$ws.Cells.Item(17, 13).ClearContents()  # M17 (was -431)
$ws.Cells.Item(17, 14).Value = -8230.143  # N17 (was -7765.8184)
$ws.Cells.Item(29, 8).Value = 1563  # H29 (was 2087.4285)
$ws.Cells.Item(29, 9).Value = 1201.75  # I29 (was 1121.2)
$ws.Cells.Item(29, 10).Value = 3008  # J29 (was 4503)
$ws.Cells.Item(29, 11).Value = 1201.75  # K29 (was 1121.2)
$ws.Cells.Item(29, 12).Value = 3008  # L29 (was 4503)
$ws.Cells.Item(29, 13).Value = -911.75  # M29 (was -831.2)
$ws.Cells.Item(29, 14).Value = -3588  # N29 (was -5083)
$ws.Cells.Item(33, 8).Value = 22009  # H33 (was 22499.334)
$ws.Cells.Item(33, 9).Value = 0  # I33 (was 0)
$ws.Cells.Item(33, 10).Value = 22009  # J33 (was 22499.334)
$ws.Cells.Item(33, 11).Value = 0  # K33 (was 0)
$ws.Cells.Item(33, 12).Value = 22009  # L33 (was 22499.334)
$ws.Cells.Item(33, 14).Value = -22513  # N33 (was -23003.334)
$ws.Cells.Item(46, 8).Value = 40015.875  # H46 (was 35592.3)
$ws.Cells.Item(46, 9).Value = 740  # I46 (was 745)
$ws.Cells.Item(46, 10).Value = 45626.715  # J46 (was 44304.125)
$ws.Cells.Item(46, 11).Value = 740  # K46 (was 745)
$ws.Cells.Item(46, 12).Value = 45626.715  # L46 (was 44304.125)
$ws.Cells.Item(46, 13).Value = -584  # M46 (was -589)
$ws.Cells.Item(46, 14).Value = -45938.715  # N46 (was -44616.125)
$ws.Cells.Item(59, 8).Value = 4000  # H59 (was 22833)
$ws.Cells.Item(59, 9).Value = 4000  # I59 (was 0)
$ws.Cells.Item(59, 10).Value = 0  # J59 (was 22833)
$ws.Cells.Item(59, 11).Value = 4000  # K59 (was 0)
$ws.Cells.Item(59, 12).ClearContents()  # L59 (was 22833)
$ws.Cells.Item(59, 13).Value = -3417  # M59 (was None)
$ws.Cells.Item(59, 14).ClearContents()  # N59 (was -23999)
$ws.Cells.Item(70, 8).Value = 8333.333000000001  # H70 (was 8500)
$ws.Cells.Item(70, 9).Value = 8333.333000000001  # I70 (was 8500)
$ws.Cells.Item(70, 10).Value = 0  # J70 (was 0)
$ws.Cells.Item(70, 11).Value = 8333.333000000001  # K70 (was 8500)
$ws.Cells.Item(70, 12).Value = 0  # L70 (was 0)
$ws.Cells.Item(70, 13).Value = -8063.333000000001  # M70 (was -8230)
$ws.Cells.Item(73, 8).Value = 8333.333000000001  # H73 (was 8500)
$ws.Cells.Item(73, 9).Value = 8333.333000000001  # I73 (was 8500)
$ws.Cells.Item(73, 10).Value = 0  # J73 (was 0)
$ws.Cells.Item(73, 11).Value = 8333.333000000001  # K73 (was 8500)
$ws.Cells.Item(73, 12).Value = 0  # L73 (was 0)
$ws.Cells.Item(73, 13).Value = -7397.333000000001  # M73 (was -7564)
$ws.Cells.Item(97, 8).Value = 1458.619  # H97 (was 989.63635)
$ws.Cells.Item(97, 9).Value = 1743.0667  # I97 (was 1287.0952)
$ws.Cells.Item(97, 10).Value = 747.5  # J97 (was 469.08334)
$ws.Cells.Item(97, 11).Value = 1743.0667  # K97 (was 1287.0952)
$ws.Cells.Item(97, 12).Value = 747.5  # L97 (was 469.08334)
$ws.Cells.Item(97, 13).Value = -1247.0667  # M97 (was -791.0952)
$ws.Cells.Item(97, 14).Value = -1739.5  # N97 (was -1461.08334)

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 3590.0312  # H22 (was 3555.5625)
$ws.Cells.Item(22, 9).Value = 2111.6428  # I22 (was 2104.7693)
$ws.Cells.Item(22, 10).Value = 4739.8887  # J22 (was 4548.2104)
$ws.Cells.Item(22, 11).Value = 2111.6428  # K22 (was 2104.7693)
$ws.Cells.Item(22, 12).Value = 4739.8887  # L22 (was 4548.2104)
$ws.Cells.Item(22, 13).Value = -1816.6428  # M22 (was -1809.7693)
$ws.Cells.Item(22, 14).Value = -5329.8887  # N22 (was -5138.2104)
$ws.Cells.Item(27, 8).Value = 3590.0312  # H27 (was 3555.5625)
$ws.Cells.Item(27, 9).Value = 2111.6428  # I27 (was 2104.7693)
$ws.Cells.Item(27, 10).Value = 4739.8887  # J27 (was 4548.2104)
$ws.Cells.Item(27, 11).Value = 2111.6428  # K27 (was 2104.7693)
$ws.Cells.Item(27, 12).Value = 4739.8887  # L27 (was 4548.2104)
$ws.Cells.Item(27, 13).Value = -2004.6428  # M27 (was -1997.7693)
$ws.Cells.Item(27, 14).Value = -4953.8887  # N27 (was -4762.2104)
$ws.Cells.Item(46, 8).Value = 4375.884  # H46 (was 4380.0713)
$ws.Cells.Item(46, 9).Value = 784.36365  # I46 (was 857)
$ws.Cells.Item(46, 10).Value = 5610.4688  # J46 (was 5630.1934)
$ws.Cells.Item(46, 11).Value = 784.36365  # K46 (was 857)
$ws.Cells.Item(46, 12).Value = 5610.4688  # L46 (was 5630.1934)
$ws.Cells.Item(46, 13).Value = -596.36365  # M46 (was -669)
$ws.Cells.Item(46, 14).Value = -5986.4688  # N46 (was -6006.1934)
$ws.Cells.Item(55, 8).Value = 1367.9524  # H55 (was 1392.1428)
$ws.Cells.Item(55, 9).Value = 1345.6666  # I55 (was 1452.625)
$ws.Cells.Item(55, 10).Value = 1384.6666  # J55 (was 1354.9231)
$ws.Cells.Item(55, 11).Value = 1345.6666  # K55 (was 1452.625)
$ws.Cells.Item(55, 12).Value = 1384.6666  # L55 (was 1354.9231)
$ws.Cells.Item(55, 13).Value = -1172.6666  # M55 (was -1279.625)
$ws.Cells.Item(55, 14).Value = -1730.6666  # N55 (was -1700.9231)
$ws.Cells.Item(100, 8).Value = 11366052  # H100 (was 12502388)
$ws.Cells.Item(100, 9).Value = 83334930  # I100 (was 125001200)
$ws.Cells.Item(100, 10).Value = 2545.7896  # J100 (was 2520.5557)
$ws.Cells.Item(100, 11).Value = 83334930  # K100 (was 125001200)
$ws.Cells.Item(100, 12).Value = 2545.7896  # L100 (was 2520.5557)
$ws.Cells.Item(100, 13).Value = -83334389  # M100 (was -125000659)
$ws.Cells.Item(100, 14).Value = -3627.7896  # N100 (was -3602.5557)
$ws.Cells.Item(132, 8).Value = 15338.8  # H132 (was 14276.556)
$ws.Cells.Item(132, 9).Value = 9997.5  # I132 (was 9997.5)
$ws.Cells.Item(132, 10).Value = 18899.666  # J132 (was 17699.8)
$ws.Cells.Item(132, 11).Value = 29992.5  # K132 (was 29992.5)
$ws.Cells.Item(132, 12).Value = 56698.99800000001  # L132 (was 53099.39999999999)
$ws.Cells.Item(132, 13).Value = -27462.5  # M132 (was -27462.5)
$ws.Cells.Item(132, 14).Value = -61758.99800000001  # N132 (was -58159.39999999999)
$ws.Cells.Item(133, 8).Value = 90109.664  # H133 (was 101664.5)
$ws.Cells.Item(133, 9).Value = 0  # I133 (was 0)
$ws.Cells.Item(133, 10).Value = 90109.664  # J133 (was 101664.5)
$ws.Cells.Item(133, 11).Value = 0  # K133 (was 0)
$ws.Cells.Item(133, 12).Value = 90109.664  # L133 (was 101664.5)
$ws.Cells.Item(133, 14).Value = -95169.664  # N133 (was -106724.5)

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 6155.6665  # H107 (was 6524.0713)
$ws.Cells.Item(107, 9).Value = 6576.1816  # I107 (was 6576.1816)
$ws.Cells.Item(107, 10).Value = 4999.25  # J107 (was 6333)
$ws.Cells.Item(107, 11).Value = 4999.25  # K107 (was 6333)
$ws.Cells.Item(107, 12).Value = 18999  # L107 (was 19728.5448)
$ws.Cells.Item(107, 13).Value = -14997.75  # M107 (was -17808.5448)
$ws.Cells.Item(107, 14).Value = -18837.75  # N107 (was -22839)
